$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking values
# like "580.63" or "132.00" are stored verbatim instead of being
# reinterpreted as numbers (which would drop trailing zeros, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

# --- Coin / Link (B, C) ---
$ws.Range("B36").Value = "Dai"
$ws.Range("B37").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# --- Price (D) ---
$ws.Range("D2").Value = "67.714.96"
$ws.Range("D3").Value = "3.332.25"
$ws.Range("D5").Value = "580.63"
$ws.Range("D6").Value = "175.67"
$ws.Range("D8").Value = "0.589"
$ws.Range("D9").Value = "3.328.96"
$ws.Range("D10").Value = "0.182"
$ws.Range("D11").Value = "0.580"
$ws.Range("D12").Value = "46.94"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D14").Value = "692.94"
$ws.Range("D15").Value = "3.875.45"
$ws.Range("D16").Value = "8.43"
$ws.Range("D17").Value = "67.703.51"
$ws.Range("D19").Value = "3.338.34"
$ws.Range("D20").Value = "17.52"
$ws.Range("D21").Value = "11.05"
$ws.Range("D22").Value = "0.893"
$ws.Range("D23").Value = "5.44"
$ws.Range("D24").Value = "16.94"
$ws.Range("D26").Value = "3.90"
$ws.Range("D27").Value = "2.68"
$ws.Range("D28").Value = "9.52"
$ws.Range("D29").Value = "32.97"
$ws.Range("D30").Value = "8.55"
$ws.Range("D31").Value = "7.07"
$ws.Range("D32").Value = "564.69"
$ws.Range("D33").Value = "10.99"
$ws.Range("D35").Value = "57.36"
$ws.Range("D36").Value = "0.998"
$ws.Range("D37").Value = "3.708.93"
$ws.Range("D38").Value = "3.31"
$ws.Range("D39").Value = "34.97"
$ws.Range("D41").Value = "3.16"
$ws.Range("D42").Value = "2.62"
$ws.Range("D43").Value = "0.0₃0671"
$ws.Range("D44").Value = "0.335"
$ws.Range("D45").Value = "3.29"
$ws.Range("D47").Value = "2.65"
$ws.Range("D48").Value = "0.128"
$ws.Range("D51").Value = "132.00"

# --- Volume(1h) (E) ---
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +6.93%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  +5.51%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("E21").Value = "  +4.45%  "
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  +5.45%  "
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  +6.24%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("E31").Value = "  +8.55%  "
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("E35").Value = "  +3.78%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  +12.52%  "
$ws.Range("E40").Value = "  +5.69%  "
$ws.Range("E41").Value = "  +7.62%  "
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("E44").Value = "  +4.80%  "
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  +6.29%  "
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("E51").Value = "  +3.36%  "

# Restore General number format on column D
$ws.Range("D2:D51").NumberFormat = "General"
